$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "MDF" component row (row 9) ---
$ws.Range("B9").Value = "MDF"
$ws.Range("C9").Value = "Madera mdf de 3 mm"
$ws.Range("D9").Value = 1

# E9 ("No aplica") gets centered like the other data cells in that column
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").Value = "No aplica"

# G9 ("SI") gets the green highlight + centered look used by the
# "Lo tenemos?" column, matching fill color FF92D050
$ws.Range("G9").Interior.Color = 5296274
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G9").Value = "SI"

# --- Remove the stray formatted cell that used to live at G12 ---
$ws.Range("G12").Clear()
$ws.Rows("12:12").AutoFit()

# --- Leave the selection where the user ended up editing ---
$ws.Range("C11").Select()
